$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 5231.1816
$ws.Range("J125").Value = 5925.5713
$ws.Range("L125").Value = 53330.14169999999
$ws.Range("N125").Value = -58250.14169999999
$ws.Range("H135").Value = 1512.8334
$ws.Range("I135").Value = 1512.8334
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 13615.5006
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -11080.5006
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 607857.9399999999
$ws.Range("I137").Value = 3303.3635
$ws.Range("J137").Value = 1119404.1
$ws.Range("K137").Value = 9910.0905
$ws.Range("L137").Value = 3358212.3
$ws.Range("M137").Value = -7360.0905
$ws.Range("N137").Value = -3363312.3
$ws.Range("H138").Value = 2022.5555
$ws.Range("I138").Value = 1641.2
$ws.Range("K138").Value = 4923.6
$ws.Range("M138").Value = 216.3999999999996

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 101138.6
$ws.Range("J23").Value = 101487.336
$ws.Range("L23").Value = 101487.336
$ws.Range("N23").Value = -102005.336
$ws.Range("H32").Value = 6316.696
$ws.Range("I32").Value = 2448.3333
$ws.Range("J32").Value = 13569.875
$ws.Range("K32").Value = 2448.3333
$ws.Range("L32").Value = 13569.875
$ws.Range("M32").Value = -2161.3333
$ws.Range("N32").Value = -14143.875
$ws.Range("H45").Value = 25005000
$ws.Range("H105").Value = 94000
$ws.Range("J105").Value = 94000
$ws.Range("L105").Value = 94000
$ws.Range("N105").Value = -100988

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2929
$ws.Range("I86").Value = 2863.0527
$ws.Range("K86").Value = 2863.0527
$ws.Range("M86").Value = -1740.0527
$ws.Range("H89").Value = 2929
$ws.Range("I89").Value = 2863.0527
$ws.Range("K89").Value = 14315.2635
$ws.Range("M89").Value = -8699.263500000001
$ws.Range("H94").Value = 4485.4287
$ws.Range("I94").Value = 3566.3333
$ws.Range("K94").Value = 3566.3333
$ws.Range("M94").Value = -3115.3333
$ws.Range("H134").Value = 4542.577
$ws.Range("I134").Value = 3112.7058
$ws.Range("J134").Value = 7243.4443
$ws.Range("K134").Value = 9338.117400000001
$ws.Range("L134").Value = 21730.3329
$ws.Range("M134").Value = -6803.117400000001
$ws.Range("N134").Value = -26800.3329

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 7120.276
$ws.Range("I7").Value = 10147.7
$ws.Range("J7").Value = 5526.8945
$ws.Range("K7").Value = 10147.7
$ws.Range("L7").Value = 5526.8945
$ws.Range("M7").Value = -10034.7
$ws.Range("N7").Value = -5752.8945
$ws.Range("H31").Value = 3351.4375
$ws.Range("I31").Value = 2326.75
$ws.Range("J31").Value = 4376.125
$ws.Range("K31").Value = 2326.75
$ws.Range("L31").Value = 4376.125
$ws.Range("M31").Value = -2031.75
$ws.Range("N31").Value = -4966.125
$ws.Range("H34").Value = 3351.4375
$ws.Range("I34").Value = 2326.75
$ws.Range("J34").Value = 4376.125
$ws.Range("K34").Value = 2326.75
$ws.Range("L34").Value = 4376.125
$ws.Range("M34").Value = -2124.75
$ws.Range("N34").Value = -4780.125
$ws.Range("H107").Value = 1670.7333
$ws.Range("I107").Value = 1196.5454
$ws.Range("J107").Value = 2974.75
$ws.Range("K107").Value = 1196.5454
$ws.Range("L107").Value = 2974.75
$ws.Range("M107").Value = 723.4546
$ws.Range("N107").Value = -6814.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 8501.5
$ws.Range("J44").Value = 12000
$ws.Range("L44").Value = 36000
$ws.Range("N44").Value = -36796
$ws.Range("H99").Value = 3750
$ws.Range("I99").Value = 3750
$ws.Range("K99").Value = 11250
$ws.Range("M99").Value = -9004
$ws.Range("H113").Value = 1870757.5
$ws.Range("I113").Value = 935.2857
$ws.Range("J113").Value = 2559639.5
$ws.Range("K113").Value = 2805.8571
$ws.Range("L113").Value = 7678918.5
$ws.Range("M113").Value = -635.8571000000002
$ws.Range("N113").Value = -7683258.5
$ws.Range("H114").Value = 5387.75
$ws.Range("I114").Value = 560.125
$ws.Range("K114").Value = 1680.375
$ws.Range("M114").Value = 1573.625
$ws.Range("H117").Value = 920.3333
$ws.Range("I117").Value = 884.6
$ws.Range("J117").Value = 1099
$ws.Range("K117").Value = 2653.8
$ws.Range("L117").Value = 3297
$ws.Range("M117").Value = 788.1999999999998
$ws.Range("N117").Value = -10181
$ws.Range("H129").Value = 55555830
$ws.Range("I129").Value = 330.4
$ws.Range("K129").Value = 991.1999999999999
$ws.Range("M129").Value = 4008.8
$ws.Range("H131").Value = 1687.8572
$ws.Range("I131").Value = 981.6667
$ws.Range("J131").Value = 2217.5
$ws.Range("K131").Value = 2945.0001
$ws.Range("L131").Value = 6652.5
$ws.Range("M131").Value = 2094.9999
$ws.Range("N131").Value = -16732.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 16333.333
$ws.Range("J58").Value = 9500
$ws.Range("L58").Value = 9500
$ws.Range("N58").Value = -10054
$ws.Range("H80").Value = 7175
$ws.Range("I80").Value = 6900
$ws.Range("J80").Value = 8000
$ws.Range("K80").Value = 6900
$ws.Range("L80").Value = 8000
$ws.Range("M80").Value = -5902
$ws.Range("N80").Value = -9996
$ws.Range("H83").Value = 7175
$ws.Range("I83").Value = 6900
$ws.Range("J83").Value = 8000
$ws.Range("K83").Value = 34500
$ws.Range("L83").Value = 40000
$ws.Range("M83").Value = -29508
$ws.Range("N83").Value = -49984
$ws.Range("H132").Value = 9216
$ws.Range("I132").Value = 13506
$ws.Range("J132").Value = 7500
$ws.Range("K132").Value = 40518
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -37988
$ws.Range("N132").Value = -27560
$ws.Range("H135").Value = 42698.91
$ws.Range("J135").Value = 42698.91
$ws.Range("L135").Value = 42698.91
$ws.Range("N135").Value = -52838.91

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8076.8076
$ws.Range("I22").Value = 1894.2778
$ws.Range("J22").Value = 21987.5
$ws.Range("K22").Value = 1894.2778
$ws.Range("L22").Value = 21987.5
$ws.Range("M22").Value = -1599.2778
$ws.Range("N22").Value = -22577.5
$ws.Range("H27").Value = 8076.8076
$ws.Range("I27").Value = 1894.2778
$ws.Range("J27").Value = 21987.5
$ws.Range("K27").Value = 1894.2778
$ws.Range("M27").Value = -1787.2778
$ws.Range("N27").Value = -22201.5
$ws.Range("H40").Value = 3475894
$ws.Range("I40").Value = 3749.6155
$ws.Range("K40").Value = 3749.6155
$ws.Range("M40").Value = -3613.6155
$ws.Range("H100").Value = 7956.5
$ws.Range("I100").Value = 8425.294
$ws.Range("K100").Value = 8425.294
$ws.Range("M100").Value = -7884.294
$ws.Range("H116").Value = 500000
$ws.Range("J116").Value = 500000
$ws.Range("L116").Value = 500000
$ws.Range("N116").Value = -509178
$ws.Range("H132").Value = 3562.8823
$ws.Range("I132").Value = 3535.5625
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 10606.6875
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -8076.6875
$ws.Range("N132").Value = -17060
$ws.Range("H136").Value = 2595.963
$ws.Range("I136").Value = 2482.2727
$ws.Range("J136").Value = 3096.2
$ws.Range("K136").Value = 7446.8181
$ws.Range("L136").Value = 9288.599999999999
$ws.Range("M136").Value = -4896.8181
$ws.Range("N136").Value = -14388.6

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 9125.666999999999
$ws.Range("I41").Value = 5000
$ws.Range("K41").Value = 5000
$ws.Range("M41").Value = -4610
$ws.Range("H54").Value = 49098.75
$ws.Range("I54").Value = 35000
$ws.Range("J54").Value = 53798.332
$ws.Range("K54").Value = 35000
$ws.Range("L54").Value = 53798.332
$ws.Range("M54").Value = -34480
$ws.Range("N54").Value = -54838.332
$ws.Range("H62").Value = 3234.1667
$ws.Range("I62").Value = 3101
$ws.Range("J62").Value = 3300.75
$ws.Range("K62").Value = 3101
$ws.Range("L62").Value = 3300.75
$ws.Range("M62").Value = -2477
$ws.Range("N62").Value = -4548.75
$ws.Range("H65").Value = 3234.1667
$ws.Range("I65").Value = 3101
$ws.Range("J65").Value = 3300.75
$ws.Range("K65").Value = 15505
$ws.Range("L65").Value = 16503.75
$ws.Range("M65").Value = -12385
$ws.Range("N65").Value = -22743.75
$ws.Range("H132").Value = 2505.3215
$ws.Range("I132").Value = 2058.6667
$ws.Range("J132").Value = 3309.3
$ws.Range("K132").Value = 6176.000100000001
$ws.Range("L132").Value = 9927.900000000001
$ws.Range("M132").Value = -3646.000100000001
$ws.Range("N132").Value = -14987.9
